# "Changes of 13th May 2022"
# Refresh the sample RTE job numbers (tracking/job/pickup/BOL ids) that are
# echoed across the RTECreation, SearchRTE, Rate and RouteDetail sheets.
#
# Each value is written as literal text (not a number) by building it as a
# quoted-string formula and then collapsing the formula to a static value
# with Copy / PasteSpecial(xlPasteValues) - this is what keeps Excel from
# reinterpreting a purely numeric-looking string (e.g. "125999607") as a
# Double and guarantees the cell keeps a shared-string ("t=s") value.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Formula = "=""" + $text + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

# --- Sheet "RTECreation" ---------------------------------------------------
$wsRTECreation = $wb.Worksheets.Item("RTECreation")
Set-TextValue $wsRTECreation "C2" "125999607"
Set-TextValue $wsRTECreation "C3" "125999629"

# --- Sheet "SearchRTE" ------------------------------------------------------
$wsSearchRTE = $wb.Worksheets.Item("SearchRTE")
Set-TextValue $wsSearchRTE "A2" "125999607"
Set-TextValue $wsSearchRTE "B2" "32395765"
Set-TextValue $wsSearchRTE "C2" "3401003"
Set-TextValue $wsSearchRTE "D2" "125999618"
Set-TextValue $wsSearchRTE "A3" "125999629"
Set-TextValue $wsSearchRTE "B3" "32395766"
Set-TextValue $wsSearchRTE "C3" "3401004"
Set-TextValue $wsSearchRTE "D3" "125999630"

# --- Sheet "Rate" ------------------------------------------------------------
$wsRate = $wb.Worksheets.Item("Rate")
Set-TextValue $wsRate "A2" "125999629"
Set-TextValue $wsRate "B2" "3401004"

# --- Sheet "RouteDetail" -----------------------------------------------------
# A2/B2 start out as blank, styled placeholder cells; once populated the
# explicit cell style is cleared (falls back to the default/Normal style).
$wsRouteDetail = $wb.Worksheets.Item("RouteDetail")
Set-TextValue $wsRouteDetail "A2" "125999629"
$wsRouteDetail.Range("A2").Style = "Normal"
Set-TextValue $wsRouteDetail "B2" "3401004"
$wsRouteDetail.Range("B2").Style = "Normal"
